# Apply the weekly "Units Completed Per Scope ID" report refresh:
#  - bump the "Report Generated On" timestamp
#  - fill in the computed Total Billed Amount
#  - clear the (now-unused) Scope ID # field
#  - populate the per-line-item Pricing column (and day TOTAL rows) that was
#    previously all zeroes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report generation timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:01 AM"

# Report summary: total billed amount now computed
$ws.Range("C8").Value = 14684.57

# Scope ID # is no longer populated
$ws.Range("G10").Value = ""

# Thursday (07/24/2025) line items + TOTAL
$ws.Range("H16").Value = 476.4
$ws.Range("H17").Value = 2858.4
$ws.Range("H18").Value = 215.84
$ws.Range("H19").Value = 3550.64

# Friday (07/25/2025) line items + TOTAL
$ws.Range("H24").Value = 238.2
$ws.Range("H25").Value = 2620.2
$ws.Range("H26").Value = 107.92
$ws.Range("H27").Value = 350.53
$ws.Range("H28").Value = 31.72
$ws.Range("H29").Value = 188.34
$ws.Range("H30").Value = 188.34
$ws.Range("H31").Value = 478.55
$ws.Range("H32").Value = 61.83
$ws.Range("H33").Value = 62.16
$ws.Range("H34").Value = 62.16
$ws.Range("H35").Value = 4389.949999999999

# Saturday (07/26/2025) line items + TOTAL
$ws.Range("H40").Value = 1191
$ws.Range("H41").Value = 833.7
$ws.Range("H42").Value = 238.2
$ws.Range("H43").Value = 238.2
$ws.Range("H44").Value = 2501.1

# Sunday (07/27/2025) line items + TOTAL
$ws.Range("H49").Value = 2858.4
$ws.Range("H50").Value = 238.2
$ws.Range("H51").Value = 215.84
$ws.Range("H52").Value = 238.2
$ws.Range("H53").Value = 476.4
$ws.Range("H54").Value = 215.84
$ws.Range("H55").Value = 4242.88
